# The author re-uploaded the workbook with 18 data rows removed from the
# "R. communis" weekly-difference table (rows whose H10_H9 value was 0 and
# were pruned during cleanup). The remaining 32 data rows keep their
# original relative order and shift up to fill rows 2:33, so the used
# range becomes A1:D33 instead of A1:D51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (Réplica, Tratamiento, H10_H9) triplets identifying each row to delete,
# matched against columns A, C, D of the current A1:D51 table.
$targets = @(
    @(6,  "C",  0),
    @(7,  "C",  0),
    @(24, "J",  0),
    @(25, "J",  0),
    @(26, "J",  0),
    @(27, "J",  0),
    @(32, "J",  0),
    @(34, "J",  0),
    @(36, "J",  0),
    @(46, "JI", 0),
    @(47, "JI", 0),
    @(49, "JI", 0),
    @(50, "JI", 0),
    @(56, "JI", 0),
    @(57, "JI", 0),
    @(76, "CI", 0),
    @(80, "CI", 0),
    @(86, "SN", 0)
)

# Walk the data rows bottom-up so deleting a matched row never shifts the
# row index of rows still to be examined above it.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
for ($r = $lastRow; $r -ge 2; $r--) {
    $a = $ws.Cells.Item($r, 1).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    $d = $ws.Cells.Item($r, 4).Value2

    foreach ($t in $targets) {
        if ($a -eq $t[0] -and $c -eq $t[1] -and $d -eq $t[2]) {
            $ws.Rows.Item($r).Delete()
            break
        }
    }
}

# Reset the sheet selection back to the default single cell (A1) so the
# saved view state no longer references the old A1:D51 extent.
$ws.Range("A1").Select()
